# Rebuild the sheet as "aperturas_masivas" style data: 8 header columns +
# 3 data rows, replacing the old COL1/COL2/COL3 demo content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$headers = @("Recaudadora", "Tipo", "Cuenta", "Fecha de otorgamiento", "Recamaras", "Banios", "Localidad", "Colonia")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Rows(1).RowHeight = 24

# ---- Data rows (rows 2-4) ----
# Columns: A Recaudadora(num) B Tipo(text) C Cuenta(num) D Fecha de otorgamiento(text)
#          E Recamaras(num) F Banios(num) G Localidad(text) H Colonia(text)

$ws.Range("A2").Value = 93
$ws.Range("B2").Value = "u"
$ws.Range("C2").Value = 123456
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "12/04/2023"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "San Agustin"
$ws.Range("H2").Value = "San Agustin"

$ws.Range("A3").Value = 77
$ws.Range("B3").Value = "u"
$ws.Range("C3").Value = 8877
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "12/04/2023"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = "El Palomar"
$ws.Range("H3").Value = "El Palomar"

$ws.Range("A4").Value = 132
$ws.Range("B4").Value = "U"
$ws.Range("C4").Value = 4455
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "12/04/2023"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = "Tlajo"
$ws.Range("H4").Value = "Tlajo"

# ---- Formatting ----
# Header row: bold font (already via existing style 1), text format + left/center alignment
$ws.Range("A1:H1").NumberFormat = "@"
$ws.Range("A1:H1").HorizontalAlignment = -4131
$ws.Range("A1:H1").VerticalAlignment = -4108

# Data rows: text format + left/center alignment (applied after values so
# numeric cells keep their numeric type, matching Excel's real behaviour)
$ws.Range("A2:H4").NumberFormat = "@"
$ws.Range("A2:H4").HorizontalAlignment = -4131
$ws.Range("A2:H4").VerticalAlignment = -4108

# ---- Column widths ----
$ws.Columns("A:C").ColumnWidth = 17.5
$ws.Columns("D").ColumnWidth = 22.1667
$ws.Columns("E:G").ColumnWidth = 17.5
$ws.Columns("H").ColumnWidth = 10.8333

# ---- Selection ----
[void]$ws.Range("B5").Select()
